$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.784.34"
$ws.Range("E2").Value = "  +0.34%  "
$ws.Range("D3").Value = "2.292.00"
$ws.Range("E3").Value = "  +0.02%  "
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("D5").Value = "'110.19"
$ws.Range("E5").Value = "  +14.35%  "
$ws.Range("D6").Value = "'268.12"
$ws.Range("E6").Value = "  -0.24%  "
$ws.Range("E7").Value = "  +0.36%  "
$ws.Range("E8").Value = "  +0.31%  "
$ws.Range("D9").Value = "'0.617"
$ws.Range("E9").Value = "  +1.10%  "
$ws.Range("D10").Value = "'47.65"
$ws.Range("E10").Value = "  +4.78%  "
$ws.Range("E11").Value = "  +1.56%  "
$ws.Range("D12").Value = "'9.05"
$ws.Range("E12").Value = "  +14.01%  "
$ws.Range("E13").Value = "  -0.08%  "
$ws.Range("D14").Value = "'15.73"
$ws.Range("E14").Value = "  +2.29%  "
$ws.Range("D15").Value = "2.631.42"
$ws.Range("E15").Value = "  -0.16%  "
$ws.Range("D16").Value = "'0.849"
$ws.Range("E16").Value = "  +0.03%  "
$ws.Range("D17").Value = "2.284.77"
$ws.Range("E17").Value = "  -0.14%  "
$ws.Range("D18").Value = "43.679.39"
$ws.Range("E18").Value = "  +0.12%  "
$ws.Range("E19").Value = "  -1.42%  "
$ws.Range("D20").Value = "'6.77"
$ws.Range("E20").Value = "  +9.18%  "
$ws.Range("D21").Value = "'72.24"
$ws.Range("E21").Value = "  +0.17%  "
$ws.Range("E22").Value = "  -4.03%  "
$ws.Range("D23").Value = "'10.03"
$ws.Range("E23").Value = "  +10.00%  "
$ws.Range("D24").Value = "'232.27"
$ws.Range("E24").Value = "  -0.24%  "
$ws.Range("D25").Value = "'2.75"
$ws.Range("E25").Value = "  +4.34%  "
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("D27").Value = "'11.64"
$ws.Range("E27").Value = "  +3.47%  "
$ws.Range("D28").Value = "'41.69"
$ws.Range("E28").Value = "  +4.67%  "
$ws.Range("E29").Value = "  -2.24%  "
$ws.Range("E30").Value = "  +2.29%  "
$ws.Range("D31").Value = "'175.75"
$ws.Range("E31").Value = "  +0.37%  "
$ws.Range("D32").Value = "'21.51"
$ws.Range("E32").Value = "  -1.57%  "
$ws.Range("D33").Value = "'0.0923"
$ws.Range("E33").Value = "  +2.92%  "
$ws.Range("D34").Value = "'5.64"
$ws.Range("E34").Value = "  +4.60%  "
$ws.Range("D35").Value = "'0.128"
$ws.Range("E35").Value = "  +1.50%  "
$ws.Range("E36").Value = "  +6.90%  "
$ws.Range("D37").Value = "'0.0364"
$ws.Range("E37").Value = "  +3.20%  "
$ws.Range("D38").Value = "'0.107"
$ws.Range("E38").Value = "  +0.19%  "
$ws.Range("E39").Value = "  +13.76%  "
$ws.Range("B40").Value = "Celestia"
$ws.Range("C40").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D40").Value = "'13.74"
$ws.Range("E40").Value = "  +11.42%  "
$ws.Range("B41").Value = "Algorand"
$ws.Range("C41").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D41").Value = "'0.242"
$ws.Range("E41").Value = "  +1.00%  "
$ws.Range("E42").Value = "  +2.12%  "
$ws.Range("D43").Value = "'73.08"
$ws.Range("E43").Value = "  +11.75%  "
$ws.Range("D44").Value = "'6.34"
$ws.Range("E44").Value = "  +23.08%  "
$ws.Range("E45").Value = "  +0.21%  "
$ws.Range("E46").Value = "  +1.81%  "
$ws.Range("D47").Value = "'8.73"
$ws.Range("E47").Value = "  -0.62%  "
$ws.Range("D48").Value = "'102.15"
$ws.Range("E48").Value = "  +4.84%  "
$ws.Range("D49").Value = "'0.0991"
$ws.Range("E49").Value = "  -2.37%  "
$ws.Range("D50").Value = "'1.23"
$ws.Range("E50").Value = "  +3.07%  "
$ws.Range("D51").Value = "'0.452"
$ws.Range("E51").Value = "  +5.58%  "
